# Signed Off Time Sheets
# As of 28/02/2014
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor Signature (initials) and signed-off date
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").NumberFormat = "mm-dd-yy"
$ws.Range("D27").Value = Get-Date -Year 2014 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Move selection to reflect the last-edited cell
$ws.Range("D27:E27").Select()
